$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.028318268037473
$ws.Cells.Item(2, 4).Value = 1.031439595793088
$ws.Cells.Item(2, 5).Value = 0.9926147277508489
$ws.Cells.Item(2, 6).Value = 1.026866817619818
$ws.Cells.Item(2, 9).Value = 1.029884976315947
$ws.Cells.Item(2, 10).Value = 1.033471251638803
$ws.Cells.Item(2, 11).Value = 1.034247709585159
$ws.Cells.Item(2, 12).Value = 0.9955398523336033
$ws.Cells.Item(2, 13).Value = 1.029688212686092
$ws.Cells.Item(2, 14).Value = 1.034938899356947
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.030407171335024
$ws.Cells.Item(3, 4).Value = 1.033021393207783
$ws.Cells.Item(3, 5).Value = 0.9936372048519304
$ws.Cells.Item(3, 6).Value = 1.029607382516777
$ws.Cells.Item(3, 9).Value = 1.030361966769969
$ws.Cells.Item(3, 10).Value = 1.035194878741642
$ws.Cells.Item(3, 11).Value = 1.035635999581812
$ws.Cells.Item(3, 12).Value = 0.9963617723202692
$ws.Cells.Item(3, 13).Value = 1.032231160345852
$ws.Cells.Item(3, 14).Value = 1.036664974207968
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.031751101258544
$ws.Cells.Item(4, 4).Value = 1.034038208046212
$ws.Cells.Item(4, 5).Value = 0.9942998659930995
$ws.Cells.Item(4, 6).Value = 1.031371866788579
$ws.Cells.Item(4, 9).Value = 1.030666365060876
$ws.Cells.Item(4, 10).Value = 1.03630240685744
$ws.Cells.Item(4, 11).Value = 1.036527144059523
$ws.Cells.Item(4, 12).Value = 0.9968940712668345
$ws.Cells.Item(4, 13).Value = 1.033867592181122
$ws.Cells.Item(4, 14).Value = 1.037774075140726
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.03231427937683
$ws.Cells.Item(5, 4).Value = 1.034464099440905
$ws.Cells.Item(5, 5).Value = 0.9945786998344013
$ws.Cells.Item(5, 6).Value = 1.032111596727522
$ws.Cells.Item(5, 9).Value = 1.030793328676452
$ws.Cells.Item(5, 10).Value = 1.036766185535673
$ws.Cells.Item(5, 11).Value = 1.036900092492043
$ws.Cells.Item(5, 12).Value = 0.9971179600051301
$ws.Cells.Item(5, 13).Value = 1.034553442003936
$ws.Cells.Item(5, 14).Value = 1.038238512437877
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.032408734490383
$ws.Cells.Item(6, 4).Value = 1.034535516828
$ws.Cells.Item(6, 5).Value = 0.9946255319796335
$ws.Cells.Item(6, 6).Value = 1.032235681558761
$ws.Cells.Item(6, 9).Value = 1.030814587789368
$ws.Cells.Item(6, 10).Value = 1.036843949937213
$ws.Cells.Item(6, 11).Value = 1.036962613963011
$ws.Cells.Item(6, 12).Value = 0.9971555583673452
$ws.Cells.Item(6, 13).Value = 1.034668477275075
$ws.Cells.Item(6, 14).Value = 1.038316387273787
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.031758633530072
$ws.Cells.Item(7, 4).Value = 1.034043904986844
$ws.Cells.Item(7, 5).Value = 0.9943035907978917
$ws.Cells.Item(7, 6).Value = 1.031381759110282
$ws.Cells.Item(7, 9).Value = 1.030668065490585
$ws.Cells.Item(7, 10).Value = 1.03630861101967
$ws.Cells.Item(7, 11).Value = 1.036532134012393
$ws.Cells.Item(7, 12).Value = 0.9968970624459043
$ws.Cells.Item(7, 13).Value = 1.033876764742984
$ws.Cells.Item(7, 14).Value = 1.037780288113577
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.029025847922547
$ws.Cells.Item(8, 4).Value = 1.031975579848342
$ws.Cells.Item(8, 5).Value = 0.9929600610674297
$ws.Cells.Item(8, 6).Value = 1.027794870959605
$ws.Cells.Item(8, 9).Value = 1.03004706351436
$ws.Cells.Item(8, 10).Value = 1.034055390335392
$ws.Cells.Item(8, 11).Value = 1.034718390870148
$ws.Cells.Item(8, 12).Value = 0.9958175282591057
$ws.Cells.Item(8, 13).Value = 1.030549515638512
$ws.Cells.Item(8, 14).Value = 1.035523867597479
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.024149278608423
$ws.Cells.Item(9, 4).Value = 1.02827815176906
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.021403931533868
$ws.Cells.Item(9, 9).Value = 1.028919746501215
$ws.Cells.Item(9, 10).Value = 1.030023798091202
$ws.Cells.Item(9, 11).Value = 1.03146612669365
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.024614870481193
$ws.Cells.Item(9, 14).Value = 1.031486550029872
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.020854566571942
$ws.Cells.Item(10, 4).Value = 1.02577579587502
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.017092234453927
$ws.Cells.Item(10, 9).Value = 1.028145269193014
$ws.Cells.Item(10, 10).Value = 1.027292714811693
$ws.Cells.Item(10, 11).Value = 1.029258326516209
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.020606776958628
$ws.Cells.Item(10, 14).Value = 1.028751588298847
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.019416970294698
$ws.Cells.Item(11, 4).Value = 1.024682941126315
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.015212260230845
$ws.Cells.Item(11, 9).Value = 1.027804308338857
$ws.Cells.Item(11, 10).Value = 1.026099331503987
$ws.Cells.Item(11, 11).Value = 1.028292512255497
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.018858177730902
$ws.Cells.Item(11, 14).Value = 1.027556510249961
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.018881284086557
$ws.Cells.Item(12, 4).Value = 1.024275568248823
$ws.Cells.Item(12, 5).Value = 0.9881042295825494
$ws.Cells.Item(12, 6).Value = 1.014511929900379
$ws.Cells.Item(12, 9).Value = 1.027676803502445
$ws.Cells.Item(12, 10).Value = 1.025654387700594
$ws.Cells.Item(12, 11).Value = 1.0279322539781
$ws.Cells.Item(12, 12).Value = 0.991905972511983
$ws.Cells.Item(12, 13).Value = 1.018206637304541
$ws.Cells.Item(12, 14).Value = 1.027110934575332
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.018996268320653
$ws.Cells.Item(13, 4).Value = 1.024363016805238
$ws.Cells.Item(13, 5).Value = 0.9881581567098647
$ws.Cells.Item(13, 6).Value = 1.014662245939842
$ws.Cells.Item(13, 9).Value = 1.027704192729392
$ws.Cells.Item(13, 10).Value = 1.025749905869864
$ws.Cells.Item(13, 11).Value = 1.028009599606739
$ws.Cells.Item(13, 12).Value = 0.9919494934313047
$ws.Cells.Item(13, 13).Value = 1.018346488060443
$ws.Cells.Item(13, 14).Value = 1.027206588391358
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.019372725266803
$ws.Cells.Item(14, 4).Value = 1.02464929711433
$ws.Cells.Item(14, 5).Value = 0.988334886381446
$ws.Cells.Item(14, 6).Value = 1.01515441243438
$ws.Cells.Item(14, 9).Value = 1.02779378631214
$ws.Cells.Item(14, 10).Value = 1.026062586600973
$ws.Cells.Item(14, 11).Value = 1.028262764247688
$ws.Cells.Item(14, 12).Value = 0.9920921077337194
$ws.Cells.Item(14, 13).Value = 1.018804363072509
$ws.Cells.Item(14, 14).Value = 1.02751971316497
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.019604446124042
$ws.Cells.Item(15, 4).Value = 1.024825492165206
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.015457381988439
$ws.Cells.Item(15, 9).Value = 1.02784887390573
$ws.Cells.Item(15, 10).Value = 1.026255017127842
$ws.Cells.Item(15, 11).Value = 1.02841854590647
$ws.Cells.Item(15, 12).Value = 0.9921799884222137
$ws.Cells.Item(15, 13).Value = 1.019086203506897
$ws.Cells.Item(15, 14).Value = 1.027712416965259
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.020949739712413
$ws.Cells.Item(16, 4).Value = 1.025848125497125
$ws.Cells.Item(16, 5).Value = 0.9890781214508735
$ws.Cells.Item(16, 6).Value = 1.017216722383279
$ws.Cells.Item(16, 9).Value = 1.02816777828623
$ws.Cells.Item(16, 10).Value = 1.027371684272969
$ws.Cells.Item(16, 11).Value = 1.029322214415237
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.020722544504028
$ws.Cells.Item(16, 14).Value = 1.028830669905816
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.021790635268107
$ws.Cells.Item(17, 4).Value = 1.026487075063732
$ws.Cells.Item(17, 5).Value = 0.9894763578477731
$ws.Cells.Item(17, 6).Value = 1.01831678307636
$ws.Cells.Item(17, 9).Value = 1.028366307347311
$ws.Cells.Item(17, 10).Value = 1.028069215766415
$ws.Cells.Item(17, 11).Value = 1.029886406821217
$ws.Cells.Item(17, 12).Value = 0.9930127773692701
$ws.Cells.Item(17, 13).Value = 1.021745430686745
$ws.Cells.Item(17, 14).Value = 1.029529191973991
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.022280061078309
$ws.Cells.Item(18, 4).Value = 1.02685886708149
$ws.Cells.Item(18, 5).Value = 0.9897087662937551
$ws.Cells.Item(18, 6).Value = 1.018957182473067
$ws.Cells.Item(18, 9).Value = 1.028481565977865
$ws.Cells.Item(18, 10).Value = 1.028475034502821
$ws.Cells.Item(18, 11).Value = 1.030214545489426
$ws.Cells.Item(18, 12).Value = 0.9932001317071766
$ws.Cells.Item(18, 13).Value = 1.022340806760161
$ws.Cells.Item(18, 14).Value = 1.02993558701955
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.02244676534472
$ws.Cells.Item(19, 4).Value = 1.026985487635342
$ws.Cells.Item(19, 5).Value = 0.9897880325774039
$ws.Cells.Item(19, 6).Value = 1.019175332751986
$ws.Cells.Item(19, 9).Value = 1.028520775027634
$ws.Cells.Item(19, 10).Value = 1.028613233208225
$ws.Cells.Item(19, 11).Value = 1.030326273053103
$ws.Cells.Item(19, 12).Value = 0.993264023964098
$ws.Cells.Item(19, 13).Value = 1.022543603632534
$ws.Cells.Item(19, 14).Value = 1.030073981982967
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.021700524559139
$ws.Cells.Item(20, 4).Value = 1.026418614819261
$ws.Cells.Item(20, 5).Value = 0.9894336180355766
$ws.Cells.Item(20, 6).Value = 1.018198886528695
$ws.Cells.Item(20, 9).Value = 1.028345063017084
$ws.Cells.Item(20, 10).Value = 1.027994485087126
$ws.Cells.Item(20, 11).Value = 1.029825972343118
$ws.Cells.Item(20, 12).Value = 0.9929783193490043
$ws.Cells.Item(20, 13).Value = 1.021635815105121
$ws.Cells.Item(20, 14).Value = 1.029454355168565
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.019261915398848
$ws.Cells.Item(21, 4).Value = 1.024565034702945
$ws.Cells.Item(21, 5).Value = 0.9882828385668255
$ws.Cells.Item(21, 6).Value = 1.015009538148483
$ws.Cells.Item(21, 9).Value = 1.02776742700398
$ws.Cells.Item(21, 10).Value = 1.025970556329891
$ws.Cells.Item(21, 11).Value = 1.028188255600064
$ws.Cells.Item(21, 12).Value = 0.9920501090198107
$ws.Cells.Item(21, 13).Value = 1.018669586974749
$ws.Cells.Item(21, 14).Value = 1.027427552200347
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.01771880905706
$ws.Cells.Item(22, 4).Value = 1.023391276255889
$ws.Cells.Item(22, 5).Value = 0.9875604150241496
$ws.Cells.Item(22, 6).Value = 1.012992519015608
$ws.Cells.Item(22, 9).Value = 1.027399279805342
$ws.Cells.Item(22, 10).Value = 1.024688358033244
$ws.Cells.Item(22, 11).Value = 1.027149792434783
$ws.Cells.Item(22, 12).Value = 0.991467000034148
$ws.Cells.Item(22, 13).Value = 1.016792804839164
$ws.Cells.Item(22, 14).Value = 1.026143533035049
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.018537790713226
$ws.Cells.Item(23, 4).Value = 1.024014311178306
$ws.Cells.Item(23, 5).Value = 0.9879432794636459
$ws.Cells.Item(23, 6).Value = 1.014062917841221
$ws.Cells.Item(23, 9).Value = 1.027594917070817
$ws.Cells.Item(23, 10).Value = 1.02536900773766
$ws.Cells.Item(23, 11).Value = 1.027701144580647
$ws.Cells.Item(23, 12).Value = 0.9917760702887607
$ws.Cells.Item(23, 13).Value = 1.017788864356329
$ws.Cells.Item(23, 14).Value = 1.026825149340116
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.021741244965208
$ws.Cells.Item(24, 4).Value = 1.026449551825202
$ws.Cells.Item(24, 5).Value = 0.9894529299347241
$ws.Cells.Item(24, 6).Value = 1.018252162754498
$ws.Cells.Item(24, 9).Value = 1.028354664083432
$ws.Cells.Item(24, 10).Value = 1.028028255874539
$ws.Cells.Item(24, 11).Value = 1.029853283001422
$ws.Cells.Item(24, 12).Value = 0.9929938892766438
$ws.Cells.Item(24, 13).Value = 1.021685349546053
$ws.Cells.Item(24, 14).Value = 1.02948817391437
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.025417489814463
$ws.Cells.Item(25, 4).Value = 1.029240474918675
$ws.Cells.Item(25, 5).Value = 0.9912096547607046
$ws.Cells.Item(25, 6).Value = 1.023064864502748
$ws.Cells.Item(25, 9).Value = 1.029215172979549
$ws.Cells.Item(25, 10).Value = 1.031073534723454
$ws.Cells.Item(25, 11).Value = 1.032313762632006
$ws.Cells.Item(25, 12).Value = 0.9944092447426411
$ws.Cells.Item(25, 13).Value = 1.026157962018926
$ws.Cells.Item(25, 14).Value = 1.032537777408549
